# Change response suggestCards Author: DEFAULT. Type: SAVE.
# Updates the RETURN column (D) of the "suggestCards" test table so that
# the card-list separator changes from ";" to "," and fixes the
# placeholder value in row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cells = @("D16", "D20", "D21", "D22", "D23", "D24")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D16").Value = "421892,450408,518761"
$ws.Range("D20").Value = "450407,421892,450408,518761"
$ws.Range("D21").Value = "450408,518841"
$ws.Range("D22").Value = "553643,459419"
$ws.Range("D23").Value = "421892,450408,518761"
$ws.Range("D24").Value = "421892"

foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "General"
}
